$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.34503388707513
$ws.Range("D2").Value = 9.025676764793712
$ws.Range("E2").Value = 14.93975322502766
$ws.Range("F2").Value = 36.79616141278254
$ws.Range("G2").Value = 3.682876781291748
$ws.Range("I2").Value = 29.7938504316008
$ws.Range("J2").Value = 11.26279986758615
$ws.Range("K2").Value = 9.231695708606876
$ws.Range("L2").Value = 10.00361220719645
$ws.Range("M2").Value = 14.81075265019728
$ws.Range("O2").Value = 27.90139393996111
$ws.Range("B3").Value = 13.22509735731198
$ws.Range("D3").Value = 9.032632520866819
$ws.Range("E3").Value = 14.97019616623903
$ws.Range("F3").Value = 36.89948319217237
$ws.Range("G3").Value = 3.684710401827016
$ws.Range("I3").Value = 29.91363246614572
$ws.Range("J3").Value = 11.27752368881742
$ws.Range("K3").Value = 8.976377409583556
$ws.Range("L3").Value = 9.982676467791443
$ws.Range("M3").Value = 14.77827274213438
$ws.Range("O3").Value = 27.98977383840431
$ws.Range("B4").Value = 13.15278558252374
$ws.Range("D4").Value = 9.037785239347059
$ws.Range("E4").Value = 14.98998817623658
$ws.Range("F4").Value = 36.97009924075528
$ws.Range("G4").Value = 3.6858968614514
$ws.Range("I4").Value = 29.99149965663736
$ws.Range("J4").Value = 11.28706726798239
$ws.Range("K4").Value = 8.816816741475469
$ws.Range("L4").Value = 9.970951468091711
$ws.Range("M4").Value = 14.75994411978592
$ws.Range("O4").Value = 28.04917860316347
$ws.Range("B5").Value = 13.12367934878537
$ws.Range("D5").Value = 9.04010732773863
$ws.Range("E5").Value = 14.99833087424754
$ws.Range("F5").Value = 37.00067824485422
$ws.Range("G5").Value = 3.686395640888251
$ws.Range("I5").Value = 30.02431926515922
$ws.Range("J5").Value = 11.29108321612628
$ws.Range("K5").Value = 8.751185826780567
$ws.Range("L5").Value = 9.966460593136214
$ws.Range("M5").Value = 14.75288581880645
$ws.Range("O5").Value = 28.07467760388186
$ws.Range("B6").Value = 13.11886888628936
$ws.Range("D6").Value = 9.040506352011256
$ws.Range("E6").Value = 14.99973294161724
$ws.Range("F6").Value = 37.00586466100337
$ws.Range("G6").Value = 3.686479387540173
$ws.Range("I6").Value = 30.02983470382716
$ws.Range("J6").Value = 11.29175773388071
$ws.Range("K6").Value = 8.740253963368977
$ws.Range("L6").Value = 9.965732317549856
$ws.Range("M6").Value = 14.75173873936509
$ws.Range("O6").Value = 28.07898963524959
$ws.Range("B7").Value = 13.1523915471253
$ws.Range("D7").Value = 9.037815654963952
$ws.Range("E7").Value = 14.99009956507709
$ws.Range("F7").Value = 36.9705043449712
$ws.Range("G7").Value = 3.685903526204432
$ws.Range("I7").Value = 29.99193786547307
$ws.Range("J7").Value = 11.28712091432849
$ws.Range("K7").Value = 8.815933955924413
$ws.Range("L7").Value = 9.970889735925194
$ws.Range("M7").Value = 14.75984725967833
$ws.Range("O7").Value = 28.04951726572109
$ws.Range("B8").Value = 13.30342017534668
$ws.Range("D8").Value = 9.027892424213613
$ws.Range("E8").Value = 14.95002208974322
$ws.Range("F8").Value = 36.8302961643949
$ws.Range("G8").Value = 3.683496462443272
$ws.Range("I8").Value = 29.83425554956582
$ws.Range("J8").Value = 11.26777243816801
$ws.Range("K8").Value = 9.144296439652463
$ws.Range("L8").Value = 9.996161021577048
$ws.Range("M8").Value = 14.79922154699387
$ws.Range("O8").Value = 27.93079999071199
$ws.Range("B9").Value = 13.60892581315117
$ws.Range("D9").Value = 9.015405040179621
$ws.Range("E9").Value = 14.88012546058818
$ws.Range("F9").Value = 36.61237524645715
$ws.Range("G9").Value = 3.679254988357498
$ws.Range("I9").Value = 29.55924651818902
$ws.Range("J9").Value = 11.23380577579654
$ws.Range("K9").Value = 9.762222151399543
$ws.Range("L9").Value = 10.05453838211056
$ws.Range("M9").Value = 14.88901935568959
$ws.Range("O9").Value = 27.73882492825158
$ws.Range("B10").Value = 13.83742581507478
$ws.Range("D10").Value = 9.010446696228975
$ws.Range("E10").Value = 14.83402837523326
$ws.Range("F10").Value = 36.48714630551928
$ws.Range("G10").Value = 3.67642765961067
$ws.Range("I10").Value = 29.37794342083683
$ws.Range("J10").Value = 11.21125152074336
$ws.Range("K10").Value = 10.19569586932578
$ws.Range("L10").Value = 10.10260516588068
$ws.Range("M10").Value = 14.96235438571532
$ws.Range("O10").Value = 27.62274054294489
$ws.Range("B11").Value = 13.94190051073367
$ws.Range("D11").Value = 9.009098767204389
$ws.Range("E11").Value = 14.81418940067027
$ws.Range("F11").Value = 36.43777016115737
$ws.Range("G11").Value = 3.675203527283136
$ws.Range("I11").Value = 29.29994662210366
$ws.Range("J11").Value = 11.20150757724405
$ws.Range("K11").Value = 10.38758036719895
$ws.Range("L11").Value = 10.1255500189486
$ws.Range("M11").Value = 14.9972433932176
$ws.Range("O11").Value = 27.57536427261901
$ws.Range("B12").Value = 13.98150870035486
$ws.Range("D12").Value = 9.008718153865088
$ws.Range("E12").Value = 14.80683878981825
$ws.Range("F12").Value = 36.42016564146635
$ws.Range("G12").Value = 3.674748852512026
$ws.Range("I12").Value = 29.27105377791849
$ws.Range("J12").Value = 11.19789165238574
$ws.Range("K12").Value = 10.45941613768696
$ws.Range("L12").Value = 10.13438960689643
$ws.Range("M12").Value = 15.01066830502535
$ws.Range("O12").Value = 27.55820616915065
$ws.Range("B13").Value = 13.97297680474622
$ws.Range("D13").Value = 9.008794363514189
$ws.Range("E13").Value = 14.80841468136598
$ws.Range("F13").Value = 36.4239084490265
$ws.Range("G13").Value = 3.674846380704651
$ws.Range("I13").Value = 29.27724779671384
$ws.Range("J13").Value = 11.19866712505312
$ws.Range("K13").Value = 10.44398274039496
$ws.Range("L13").Value = 10.13247920348335
$ws.Range("M13").Value = 15.00776763524702
$ws.Range("O13").Value = 27.56186666290809
$ws.Range("B14").Value = 13.94515835254834
$ws.Range("D14").Value = 9.009064856277318
$ws.Range("E14").Value = 14.81358141880592
$ws.Range("F14").Value = 36.43629991016073
$ws.Range("G14").Value = 3.675165943239624
$ws.Range("I14").Value = 29.29755671351101
$ws.Range("J14").Value = 11.20120861378118
$ws.Range("K14").Value = 10.39350725228996
$ws.Range("L14").Value = 10.12627425864771
$ws.Range("M14").Value = 14.9983436344829
$ws.Range("O14").Value = 27.57393698038613
$ws.Range("B15").Value = 13.92812383233789
$ws.Range("D15").Value = 9.009247425797293
$ws.Range("E15").Value = 14.81676727036166
$ws.Range("F15").Value = 36.44403244128813
$ws.Range("G15").Value = 3.675362839378638
$ws.Range("I15").Value = 29.31008019444708
$ws.Range("J15").Value = 11.20277496290187
$ws.Range("K15").Value = 10.36248010277944
$ws.Range("L15").Value = 10.1224930684489
$ws.Range("M15").Value = 14.99259872749077
$ws.Range("O15").Value = 27.58143230854273
$ws.Range("B16").Value = 13.83060620720842
$ws.Range("D16").Value = 9.010552991501802
$ws.Range("E16").Value = 14.83534759968458
$ws.Range("F16").Value = 36.49052604408271
$ws.Range("G16").Value = 3.676508904147432
$ws.Range("I16").Value = 29.38313071756188
$ws.Range("J16").Value = 11.21189866823922
$ws.Range("K16").Value = 10.18304341651643
$ws.Range("L16").Value = 10.1011270151984
$ws.Range("M16").Value = 14.96010450592863
$ws.Range("O16").Value = 27.625946092063
$ws.Range("B17").Value = 13.77089623953396
$ws.Range("D17").Value = 9.011585894379483
$ws.Range("E17").Value = 14.84703520831471
$ws.Range("F17").Value = 36.52099381494385
$ws.Range("G17").Value = 3.677227835516067
$ws.Range("I17").Value = 29.42909112583853
$ws.Range("J17").Value = 11.21762772326681
$ws.Range("K17").Value = 10.07155869688465
$ws.Range("L17").Value = 10.08829299460348
$ws.Range("M17").Value = 14.9405572227253
$ws.Range("O17").Value = 27.65464582318095
$ws.Range("B18").Value = 13.73660422337798
$ws.Range("D18").Value = 9.012265468461274
$ws.Range("E18").Value = 14.85386408732163
$ws.Range("F18").Value = 36.53923250468159
$ws.Range("G18").Value = 3.677647187264493
$ws.Range("I18").Value = 29.45594798098218
$ws.Range("J18").Value = 11.22097152089862
$ws.Range("K18").Value = 10.00693973550747
$ws.Range("L18").Value = 10.08101303604266
$ws.Range("M18").Value = 14.92945842458094
$ws.Range("O18").Value = 27.67166427842325
$ws.Range("B19").Value = 13.7250033170629
$ws.Range("D19").Value = 9.012510261969227
$ws.Range("E19").Value = 14.85619453507535
$ws.Range("F19").Value = 36.54553046513438
$ws.Range("G19").Value = 3.677790177144392
$ws.Range("I19").Value = 29.46511370591635
$ws.Range("J19").Value = 11.22211202963954
$ws.Range("K19").Value = 9.984977742106748
$ws.Range("L19").Value = 10.07856578492933
$ws.Range("M19").Value = 14.92572554029424
$ws.Range("O19").Value = 27.6775141817427
$ws.Range("B20").Value = 13.77724734861828
$ws.Range("D20").Value = 9.011467098687845
$ws.Range("E20").Value = 14.84578002720121
$ws.Range("F20").Value = 36.51767651548095
$ws.Range("G20").Value = 3.677150699788265
$ws.Range("I20").Value = 29.42415493222818
$ws.Range("J20").Value = 11.21701282832794
$ws.Range("K20").Value = 10.08347824161478
$ws.Range("L20").Value = 10.08964868952302
$ws.Range("M20").Value = 14.94262317883974
$ws.Range("O20").Value = 27.65153777837853
$ws.Range("B21").Value = 13.95332830004562
$ws.Range("D21").Value = 9.008981888489625
$ws.Range("E21").Value = 14.81205943174296
$ws.Range("F21").Value = 36.43263055845455
$ws.Range("G21").Value = 3.675071839381294
$ws.Range("I21").Value = 29.29157405499163
$ws.Range("J21").Value = 11.20046011397506
$ws.Range("K21").Value = 10.40835604262236
$ws.Range("L21").Value = 10.12809274176572
$ws.Range("M21").Value = 15.00110595938435
$ws.Range("O21").Value = 27.57037039539492
$ws.Range("B22").Value = 14.06866106827364
$ws.Range("D22").Value = 9.008113994661317
$ws.Range("E22").Value = 14.79096491738746
$ws.Range("F22").Value = 36.3834204111503
$ws.Range("G22").Value = 3.673764909626503
$ws.Range("I22").Value = 29.20867128228449
$ws.Range("J22").Value = 11.19007252924834
$ws.Range("K22").Value = 10.61583830329998
$ws.Range("L22").Value = 10.1540954400318
$ws.Range("M22").Value = 15.04056732887578
$ws.Range("O22").Value = 27.52188297711038
$ws.Range("B23").Value = 14.00709248337312
$ws.Range("D23").Value = 9.008508239851791
$ws.Range("E23").Value = 14.8021373022464
$ws.Range("F23").Value = 36.40910127104191
$ws.Range("G23").Value = 3.674457723774538
$ws.Range("I23").Value = 29.25257563625451
$ws.Range("J23").Value = 11.19557728619464
$ws.Range("K23").Value = 10.50556389771862
$ws.Range("L23").Value = 10.14013849647362
$ws.Range("M23").Value = 15.01939487554265
$ws.Range("O23").Value = 27.54734396812293
$ws.Range("B24").Value = 13.774375897644
$ws.Range("D24").Value = 9.011520539098036
$ws.Range("E24").Value = 14.84634715341165
$ws.Range("F24").Value = 36.51917401657509
$ws.Range("G24").Value = 3.677185554069255
$ws.Range("I24").Value = 29.4263852346113
$ws.Range("J24").Value = 11.21729066628989
$ws.Range("K24").Value = 10.07809104550581
$ws.Range("L24").Value = 10.08903547273668
$ws.Range("M24").Value = 14.94168872600412
$ws.Range("O24").Value = 27.65294131012129
$ws.Range("B25").Value = 13.52544783836396
$ws.Range("D25").Value = 9.018040232544298
$ws.Range("E25").Value = 14.89810813045954
$ws.Range("F25").Value = 36.66521113875774
$ws.Range("G25").Value = 3.680351474113276
$ws.Range("I25").Value = 29.62999297546332
$ws.Range("J25").Value = 11.24257138652865
$ws.Range("K25").Value = 9.598346821690253
$ws.Range("L25").Value = 10.0378214437589
$ws.Range("M25").Value = 14.86341052512986
$ws.Range("O25").Value = 27.78638085103452
